# Scheduled runner update: refresh market-price-derived profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across several
# sheets, per the latest pull of Asura market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 552.6316
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 552.6316
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1657.8948
$ws.Range("N17").Value = -1993.8948

$ws.Range("H76").Value = 6100.375
$ws.Range("I76").Value = 6100.375
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6100.375
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -5785.375
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 6100.375
$ws.Range("I79").Value = 6100.375
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6100.375
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -5008.375
$ws.Range("N79").ClearContents()

$ws.Range("H126").Value = 44983.332
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 44983.332
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 44983.332
$ws.Range("N126").Value = -54863.332

$ws.Range("H137").Value = 1156.6216
$ws.Range("I137").Value = 1126.3939
$ws.Range("J137").Value = 1406
$ws.Range("K137").Value = 3379.1817
$ws.Range("L137").Value = 4218
$ws.Range("M137").Value = -829.1817000000001
$ws.Range("N137").Value = -9318

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2232
$ws.Range("M4").ClearContents()

$ws.Range("H9").Value = 50000
$ws.Range("I9").Value = 50000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 50000
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -49830
$ws.Range("N9").ClearContents()

$ws.Range("H13").Value = 49900
$ws.Range("I13").Value = 50000
$ws.Range("J13").Value = 49800
$ws.Range("K13").Value = 50000
$ws.Range("L13").Value = 49800
$ws.Range("M13").Value = -49856
$ws.Range("N13").Value = -50088

$ws.Range("H20").Value = 50000
$ws.Range("I20").Value = 50000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 50000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -49730
$ws.Range("N20").ClearContents()

$ws.Range("H123").Value = 24128.875
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24128.875
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24128.875
$ws.Range("N123").Value = -33928.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 48000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 48000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 48000
$ws.Range("N126").Value = -57880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1934.4
$ws.Range("I31").Value = 1624.3478
$ws.Range("J31").Value = 5500
$ws.Range("K31").Value = 1624.3478
$ws.Range("L31").Value = 5500
$ws.Range("M31").Value = -1329.3478
$ws.Range("N31").Value = -6090

$ws.Range("H34").Value = 1934.4
$ws.Range("I34").Value = 1624.3478
$ws.Range("J34").Value = 5500
$ws.Range("K34").Value = 1624.3478
$ws.Range("L34").Value = 5500
$ws.Range("M34").Value = -1422.3478
$ws.Range("N34").Value = -5904

$ws.Range("H132").Value = 713535.6
$ws.Range("I132").Value = 902578.6
$ws.Range("J132").Value = 4624.5
$ws.Range("K132").Value = 2707735.8
$ws.Range("L132").Value = 13873.5
$ws.Range("M132").Value = -2705205.8
$ws.Range("N132").Value = -18933.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 842.7857
$ws.Range("I2").Value = 1519.8572
$ws.Range("J2").Value = 165.71428
$ws.Range("K2").Value = 9119.143199999999
$ws.Range("L2").Value = 994.28568
$ws.Range("M2").Value = -9006.143199999999
$ws.Range("N2").Value = -1220.28568

$ws.Range("H17").Value = 5025
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 10000
$ws.Range("K17").Value = 150
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 19
$ws.Range("N17").Value = -30338

$ws.Range("H34").Value = 508.81482
$ws.Range("I34").Value = 274.33334
$ws.Range("J34").Value = 696.4
$ws.Range("K34").Value = 823.0000200000001
$ws.Range("L34").Value = 2089.2
$ws.Range("M34").Value = -739.0000200000001
$ws.Range("N34").Value = -2257.2

$ws.Range("H39").Value = 3369.75
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 4326.3335
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 12979.0005
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -13567.0005

$ws.Range("H55").Value = 3685.5715
$ws.Range("I55").Value = 1800.4
$ws.Range("J55").Value = 8398.5
$ws.Range("K55").Value = 5401.200000000001
$ws.Range("L55").Value = 25195.5
$ws.Range("M55").Value = -5224.200000000001
$ws.Range("N55").Value = -25549.5

$ws.Range("H109").Value = 2402
$ws.Range("I109").Value = 700
$ws.Range("J109").Value = 3536.6667
$ws.Range("K109").Value = 2100
$ws.Range("L109").Value = 10610.0001
$ws.Range("M109").Value = -1060
$ws.Range("N109").Value = -12690.0001

$ws.Range("H131").Value = 30304390
$ws.Range("I131").Value = 430
$ws.Range("J131").Value = 45456372
$ws.Range("K131").Value = 1290
$ws.Range("L131").Value = 136369116
$ws.Range("M131").Value = 3750
$ws.Range("N131").Value = -136379196

$ws.Range("H136").Value = 5176.125
$ws.Range("I136").Value = 1074.8334
$ws.Range("J136").Value = 7636.9
$ws.Range("K136").Value = 3224.5002
$ws.Range("L136").Value = 22910.7
$ws.Range("M136").Value = 1875.4998
$ws.Range("N136").Value = -33110.7

$ws.Range("H137").Value = 1857.1666
$ws.Range("I137").Value = 1257.7778
$ws.Range("J137").Value = 3655.3333
$ws.Range("K137").Value = 3773.3334
$ws.Range("L137").Value = 10965.9999
$ws.Range("M137").Value = 1326.6666
$ws.Range("N137").Value = -21165.9999

$ws.Range("H138").Value = 2580.889
$ws.Range("I138").Value = 931.25
$ws.Range("J138").Value = 3900.6
$ws.Range("K138").Value = 2793.75
$ws.Range("L138").Value = 11701.8
$ws.Range("M138").Value = 2346.25
$ws.Range("N138").Value = -21981.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 6000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 6000
$ws.Range("N12").Value = -6280

$ws.Range("H52").Value = 23333.334
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 23333.334
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 23333.334
$ws.Range("N52").Value = -23851.334
$ws.Range("M52").ClearContents()

$ws.Range("H70").Value = 5827.136
$ws.Range("I70").Value = 5733.1333
$ws.Range("J70").Value = 6028.5713
$ws.Range("K70").Value = 5733.1333
$ws.Range("L70").Value = 6028.5713
$ws.Range("M70").Value = -5463.1333
$ws.Range("N70").Value = -6568.5713

$ws.Range("H73").Value = 5827.136
$ws.Range("I73").Value = 5733.1333
$ws.Range("J73").Value = 6028.5713
$ws.Range("K73").Value = 5733.1333
$ws.Range("L73").Value = 6028.5713
$ws.Range("M73").Value = -4797.1333
$ws.Range("N73").Value = -7900.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 24013.666
$ws.Range("I45").Value = 26020.5
$ws.Range("J45").Value = 20000
$ws.Range("K45").Value = 26020.5
$ws.Range("L45").Value = 20000
$ws.Range("M45").Value = -25613.5
$ws.Range("N45").Value = -20814

$ws.Range("H48").Value = 25000
$ws.Range("I48").Value = 50000
$ws.Range("J48").Value = 16666.666
$ws.Range("K48").Value = 50000
$ws.Range("L48").Value = 16666.666
$ws.Range("M48").Value = -49339
$ws.Range("N48").Value = -17988.666

$ws.Range("H132").Value = 4368.925
$ws.Range("I132").Value = 4588.069
$ws.Range("J132").Value = 3791.182
$ws.Range("K132").Value = 13764.207
$ws.Range("L132").Value = 11373.546
$ws.Range("M132").Value = -11234.207
$ws.Range("N132").Value = -16433.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 15000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 15000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15826

$ws.Range("H42").Value = 32511
$ws.Range("I42").Value = 50022
$ws.Range("J42").Value = 15000
$ws.Range("K42").Value = 50022
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = -49644
$ws.Range("N42").Value = -15756

$ws.Range("H43").Value = 35000
$ws.Range("I43").Value = 50000
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 50000
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = -49851
$ws.Range("N43").Value = -20298

$ws.Range("H113").Value = 227.33333
$ws.Range("I113").Value = 242.22223
$ws.Range("J113").Value = 182.66667
$ws.Range("K113").Value = 726.66669
$ws.Range("L113").Value = 548.00001
$ws.Range("M113").Value = 1443.33331
$ws.Range("N113").Value = -4888.00001

$ws.Range("H136").Value = 1600
$ws.Range("I136").Value = 1482.7587
$ws.Range("J136").Value = 2166.6667
$ws.Range("K136").Value = 4448.2761
$ws.Range("L136").Value = 6500.000100000001
$ws.Range("M136").Value = -1898.2761
$ws.Range("N136").Value = -11600.0001
